$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H1 - copy the formatting of the
# neighboring header cell (G1) so the new header matches the existing
# header style, then overwrite the copied text with "Save".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the Save values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
